$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("desguace")
$ws.Range("A5").Value = "desguace"
$ws.Range("B5").Value = "E/P. CAPRICORNIO 7"
$ws.Range("C5").Value = "Embarcación"
$ws.Range("D5").Value = "AS/42"
$ws.Range("E5").Value = "AS/42-121"
